$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..5) {
    $ws.Cells.Item($r, 3).Value = 45184
}
